$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 298-299; everything from old row 298 onward
# shifts down by two (old 298->300, ..., old 398->400).
$ws.Rows("298:299").Insert()

# Fill the two newly inserted rows with the new "Early Glo" records.
$ws.Range("A298").Value = 9
$ws.Range("B298").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C298").Value = "Metropolitana"
$ws.Range("D298").Value = 44524
$ws.Range("E298").Value = 13
$ws.Range("F298").Value = "Fruta"
$ws.Range("G298").Value = 100103
$ws.Range("H298").Value = "Frutos de hueso (carozo)"
$ws.Range("I298").Value = 100103006
$ws.Range("J298").Value = "Nectarín"
$ws.Range("K298").Value = "Early Glo"
$ws.Range("L298").Value = "Primera"
$ws.Range("M298").Value = 300
$ws.Range("N298").Value = 10000
$ws.Range("O298").Value = 10000
$ws.Range("P298").Value = 10000
$ws.Range("Q298").Value = '$/bandeja 8 kilos empedrada'
$ws.Range("R298").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S298").Value = 1250
$ws.Range("T298").Value = 8

$ws.Range("A299").Value = 9
$ws.Range("B299").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C299").Value = "Metropolitana"
$ws.Range("D299").Value = 44524
$ws.Range("E299").Value = 13
$ws.Range("F299").Value = "Fruta"
$ws.Range("G299").Value = 100103
$ws.Range("H299").Value = "Frutos de hueso (carozo)"
$ws.Range("I299").Value = 100103006
$ws.Range("J299").Value = "Nectarín"
$ws.Range("K299").Value = "Early Glo"
$ws.Range("L299").Value = "Segunda"
$ws.Range("M299").Value = 350
$ws.Range("N299").Value = 8000
$ws.Range("O299").Value = 8000
$ws.Range("P299").Value = 8000
$ws.Range("Q299").Value = '$/bandeja 8 kilos empedrada'
$ws.Range("R299").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S299").Value = 1000
$ws.Range("T299").Value = 8
